$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns keep their textual representation (values that parse as plain
# numbers would otherwise be auto-converted to numbers and lose formatting,
# e.g. trailing zeros or switch to scientific notation).

$ws.Range("D2").Value = "27.215.15"
$ws.Range("E2").Value = "  +0.85%  "

$ws.Range("D3").Value = "1.853.09"
$ws.Range("E3").Value = "  +1.54%  "

$ws.Range("E4").Value = "  -0.48%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.16"
$ws.Range("E5").Value = "  +0.47%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4634"
$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3724"
$ws.Range("E8").Value = "  +0.59%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07285"
$ws.Range("E9").Value = "  -0.75%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8877"
$ws.Range("E10").Value = "  +1.46%  "

$ws.Range("E11").Value = "  +1.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07816"
$ws.Range("E12").Value = "  -1.24%  "

$ws.Range("D13").Value = "1.918.80"
$ws.Range("E13").Value = "  +5.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.377"
$ws.Range("E14").Value = "  +0.79%  "

$ws.Range("E15").Value = "  -0.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.13"
$ws.Range("E16").Value = "  -0.14%  "

$ws.Range("E17").Value = "  -0.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008912"
$ws.Range("E18").Value = "  +0.52%  "

$ws.Range("E19").Value = "  -0.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.71"
$ws.Range("E20").Value = "  -0.33%  "

$ws.Range("D21").Value = "27.237.42"
$ws.Range("E21").Value = "  +0.81%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.061"
$ws.Range("E22").Value = "  -0.81%  "

$ws.Range("E23").Value = "  -0.37%  "

$ws.Range("D24").Value = "2.143.71"
$ws.Range("E24").Value = "  +4.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.949"
$ws.Range("E25").Value = "  +5.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.76"
$ws.Range("E26").Value = "  -0.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.41"
$ws.Range("E27").Value = "  -0.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.041"
$ws.Range("E28").Value = "  +0.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.74"
$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.054"
$ws.Range("E30").Value = "  -1.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08820"
$ws.Range("E31").Value = "  -0.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.186"
$ws.Range("E32").Value = "  +7.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7663"
$ws.Range("E33").Value = "  +5.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.169"
$ws.Range("E34").Value = "  +3.40%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.504"
$ws.Range("E35").Value = "  +1.49%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.737"
$ws.Range("E36").Value = "  +10.84%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.090"
$ws.Range("E37").Value = "  +1.86%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01940"
$ws.Range("E38").Value = "  -0.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05229"
$ws.Range("E39").Value = "  +0.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.933"
$ws.Range("E40").Value = "  -0.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.057"
$ws.Range("E41").Value = "  -0.69%  "

$ws.Range("E42").Value = "  -1.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1627"
$ws.Range("E43").Value = "  +0.31%  "

$ws.Range("E44").Value = "  +2.74%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4788"
$ws.Range("E45").Value = "  -0.92%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.29"
$ws.Range("E46").Value = "  +1.18%  "

$ws.Range("E47").Value = "  -0.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.07"
$ws.Range("E48").Value = "  -0.68%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.637"
$ws.Range("E49").Value = "  +0.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06200"
$ws.Range("E50").Value = "  +0.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.58"
$ws.Range("E51").Value = "  +1.14%  "
